# Update AdaBoostRegressor result rows across all 8 sheets with the
# "Random" search type, the winning hyperparameter dict, and the five
# metric values (train/test side in columns B:H, and the second
# train/test side in columns J:P).
#
# NOTE: this interpreter only reliably binds *positional* function
# parameters (named "-param value" binding does not populate the
# parameter variables), so Set-AdaRow below is called positionally.

$wb = $excel.ActiveWorkbook

function Set-AdaRow {
    param(
        $ws,
        $row,
        $bVal,
        $cVal,
        $dVal,
        $eVal,
        $fVal,
        $gVal,
        $hVal,
        $hasRight,
        $jVal,
        $kVal,
        $lVal,
        $mVal,
        $nVal,
        $oVal,
        $pVal
    )

    $ws.Cells.Item($row, 2).Value2 = $bVal
    $ws.Cells.Item($row, 3).Value2 = $cVal
    $ws.Cells.Item($row, 4).Value2 = $dVal
    $ws.Cells.Item($row, 5).Value2 = $eVal
    $ws.Cells.Item($row, 6).Value2 = $fVal
    $ws.Cells.Item($row, 7).Value2 = $gVal
    $ws.Cells.Item($row, 8).Value2 = $hVal

    if ($hasRight) {
        $ws.Cells.Item($row, 10).Value2 = $jVal
        $ws.Cells.Item($row, 11).Value2 = $kVal
        $ws.Cells.Item($row, 12).Value2 = $lVal
        $ws.Cells.Item($row, 13).Value2 = $mVal
        $ws.Cells.Item($row, 14).Value2 = $nVal
        $ws.Cells.Item($row, 15).Value2 = $oVal
        $ws.Cells.Item($row, 16).Value2 = $pVal
    }
}

# Sheet 1: "Option 1 - LR1 - DN1 (70-30)" -- AdaBoostRegressor is row 17
$ws1 = $wb.Worksheets.Item(1)
Set-AdaRow $ws1 17 `
    "Random" "{'n_estimators': 50, 'loss': 'linear', 'learning_rate': 0.01}" `
    0.05585056819295494 0.1616071064464095 2.08185743933642 0.236327248096691 33.75796421693686 `
    $true `
    "Random" "{'n_estimators': 400, 'loss': 'square', 'learning_rate': 1.0}" `
    0.07980255259731403 0.1958961814587728 2.420045338215218 0.2824934558486515 38.19700835405555

# Sheet 2: "Option 1 - LR1 - DN2 (70-30)" -- AdaBoostRegressor is row 18
$ws2 = $wb.Worksheets.Item(2)
Set-AdaRow $ws2 18 `
    "Random" "{'n_estimators': 400, 'loss': 'square', 'learning_rate': 1.0}" `
    0.9729227209178051 0.5187021153095465 0.5143165607292146 0.9863684508933794 25.19054015180987 `
    $true `
    "Random" "{'n_estimators': 400, 'loss': 'square', 'learning_rate': 1.0}" `
    0.9728612321050397 0.5258878610193886 0.5155716272343706 0.9863372811087695 25.46254193707204

# Sheet 3: "Option 1 - LR2 - DN1 (70-30)" -- AdaBoostRegressor is row 17
$ws3 = $wb.Worksheets.Item(3)
Set-AdaRow $ws3 17 `
    "Random" "{'n_estimators': 50, 'loss': 'exponential', 'learning_rate': 0.01}" `
    36.94739017274242 4.710618455963436 1.958450970928508 6.078436490804393 40.29281576557526 `
    $true `
    "Random" "{'n_estimators': 50, 'loss': 'exponential', 'learning_rate': 0.01}" `
    44.14785569163372 5.163666951529017 1.997047517039016 6.64438527567703 43.64314446519597

# Sheet 4: "Option 1 - LR2 - DN2 (70-30)" -- AdaBoostRegressor is row 17
$ws4 = $wb.Worksheets.Item(4)
Set-AdaRow $ws4 17 `
    "Random" "{'n_estimators': 50, 'loss': 'exponential', 'learning_rate': 0.01}" `
    41.55191644077423 4.955605281593018 1.992177650435151 6.446077601206351 42.78056161910703 `
    $true `
    "Random" "{'n_estimators': 50, 'loss': 'linear', 'learning_rate': 0.01}" `
    37.25401664781091 4.925836155158666 1.715073773233124 6.103606855606848 43.05179853034294

# Sheet 5: "Option 1 - NLR1 - DN1 (70-30)" -- AdaBoostRegressor is row 17
$ws5 = $wb.Worksheets.Item(5)
Set-AdaRow $ws5 17 `
    "Random" "{'n_estimators': 400, 'loss': 'linear', 'learning_rate': 1.0}" `
    0.06135704708703064 0.1764948092588685 24957820681211.85 0.2477035467792713 37.32149410301708 `
    $true `
    "Random" "{'n_estimators': 400, 'loss': 'linear', 'learning_rate': 1.0}" `
    0.06135704708703064 0.1764948092588685 24957820681211.85 0.2477035467792713 37.32149410301708

# Sheet 6: "Option 1 - NLR1 - DN2 (70-30)" -- AdaBoostRegressor is row 17
$ws6 = $wb.Worksheets.Item(6)
Set-AdaRow $ws6 17 `
    "Random" "{'n_estimators': 400, 'loss': 'square', 'learning_rate': 1.0}" `
    1.168131588510747 0.556534856715671 0.6792500241422448 1.08080136404001 25.5579733083827 `
    $true `
    "Random" "{'n_estimators': 400, 'loss': 'square', 'learning_rate': 1.0}" `
    1.168131588510747 0.556534856715671 0.6792500241422448 1.08080136404001 25.5579733083827

# Sheet 7: "Option 1 - NLR2 - DN1 (70-30)" -- AdaBoostRegressor is row 17
# Only the left-hand block (B17:H17) is updated for this sheet; J17:P17
# remain untouched/empty.
$ws7 = $wb.Worksheets.Item(7)
Set-AdaRow $ws7 17 `
    "Random" "{'n_estimators': 400, 'loss': 'square', 'learning_rate': 1.0}" `
    0.07906747106674768 0.203035104995928 2.493391640330773 0.2811893864759971 38.71815714166133 `
    $false `
    $null $null $null $null $null $null $null

# Sheet 8: "Option 1 - NLR2 - DN2 (70-30)" -- AdaBoostRegressor is row 17
# Only the left-hand block (B17:H17) is updated for this sheet; J17:P17
# remain untouched/empty.
$ws8 = $wb.Worksheets.Item(8)
Set-AdaRow $ws8 17 `
    "Random" "{'n_estimators': 400, 'loss': 'square', 'learning_rate': 1.0}" `
    0.8701398583310617 0.4987169437303953 0.4654834554096153 0.9328128742309798 27.41888005021302 `
    $false `
    $null $null $null $null $null $null $null
